# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column O ("N:N") on the
#   "Repayment Schedule" sheet, pushing the old "Late" header/column and the
#   "Outstanding" column one position to the right.
# - Make "Repayment Schedule" the active sheet (it was "Summary" before),
#   and leave the selection on cell R9 (matches the post-edit selection).

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment Schedule")

# Switch to the Repayment Schedule tab first (matches workbookView
# activeTab going from the Summary sheet to this one, and tabSelected
# moving off Summary onto Repayment Schedule).
$wsRepayment.Activate()

# Insert a blank column at N, shifting existing N..P columns (Late,
# Outstanding, etc.) one column to the right (N->O, P->Q).
$wsRepayment.Columns("N:N").Insert()

# Final selection left on the sheet after the edit.
$wsRepayment.Range("R9").Select()
